$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Cell value edits -------------------------------------------------
# Row 3: email (D3) removed entirely, leaving only the (hyperlink) style behind.
$ws.Range("D3").ClearContents()

# Row 4: the stray mobile number in E4 is removed.
$ws.Range("E4").ClearContents()

# Row 6: email updated to a new (longer) address + mobile number changed.
# (Set before the row-2 edit below so new shared strings land in the same
# order as the target workbook.)
$ws.Range("D6").Value = "afs231443453453453453452222222ar1@gmail.com"
$ws.Range("E6").Value = 3532342369342230

# Row 2: email + mobile number replaced with new test data.
$ws.Range("D2").Value = "eqewwewrq@dsa.com"
$ws.Range("E2").Value = 99

# --- Hyperlinks ---------------------------------------------------------
# The engine's Hyperlinks.Delete() call clears every hyperlink on the
# sheet, so wipe them all and re-add the ones that remain (skipping D3,
# whose email/hyperlink was deleted above) in the target relationship
# order: D5, D6, D7, D8, D9, D4, D10, D11, D2.
# Hyperlinks.Add() re-applies the "Hyperlink" cell style but as a *new*
# style entry, so explicitly restoring the named style afterwards keeps
# every hyperlinked cell pointed at the original style index.
$ws.Cells.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:saquibrah321mani007@gmail.com")
$ws.Range("D5").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:afs231443453453453453452222222ar1@gmail.com")
$ws.Range("D6").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("D7"), "mailto:afs222ar@gmail.com")
$ws.Range("D7").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("D8"), "mailto:a2f222456sdsar@gmail.com")
$ws.Range("D8").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("D9"), "mailto:a2222222222fsar@gmail.com")
$ws.Range("D9").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:cool12333@gmail.com")
$ws.Range("D4").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("D10"), "mailto:geddfsar@gmail.com")
$ws.Range("D10").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("D11"), "mailto:geasdfasdfddfsar@gmail.com")
$ws.Range("D11").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:eqewwewrq@dsa.com")
$ws.Range("D2").Style = "Hyperlink"

# --- View state -----------------------------------------------------
# Selection moves from X3 (with the view scrolled to show column M) back
# to F2 with the view scrolled to the default top-left (A1).
$ws.Range("F2").Select()
